$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.047.02"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "2.048.95"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'248.23"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "'0.662"
$ws.Range("E6").Value = "  +1.25%  "

$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'57.70"
$ws.Range("E7").Value = "  +5.16%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.379"
$ws.Range("E9").Value = "  +0.83%  "

$ws.Range("D10").Value = "'0.0776"
$ws.Range("E10").Value = "  -1.41%  "

$ws.Range("E11").Value = "  +1.56%  "

$ws.Range("D12").Value = "'15.71"
$ws.Range("E12").Value = "  +4.10%  "

$ws.Range("D13").Value = "2.345.88"
$ws.Range("E13").Value = "  -0.35%  "

$ws.Range("D14").Value = "'0.807"
$ws.Range("E14").Value = "  -1.05%  "

$ws.Range("D15").Value = "'5.51"
$ws.Range("E15").Value = "  +5.59%  "

$ws.Range("D16").Value = "2.046.00"
$ws.Range("E16").Value = "  -0.53%  "

$ws.Range("D17").Value = "37.035.08"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").Value = "'16.64"
$ws.Range("E18").Value = "  +16.63%  "

$ws.Range("D19").Value = "'74.55"
$ws.Range("E19").Value = "  +3.15%  "

$ws.Range("D20").Value = "0.0₃0896"
$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").Value = "'5.32"
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").Value = "'235.82"
$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "'2.37"
$ws.Range("E24").Value = "  -2.22%  "

$ws.Range("E25").Value = "  +11.16%  "

$ws.Range("D26").Value = "'167.50"
$ws.Range("E26").Value = "  -1.55%  "

$ws.Range("D27").Value = "'9.13"
$ws.Range("E27").Value = "  +0.82%  "

$ws.Range("D28").Value = "'19.71"
$ws.Range("E28").Value = "  -2.47%  "

$ws.Range("E30").Value = "  +5.12%  "

$ws.Range("D31").Value = "'4.67"
$ws.Range("E31").Value = "  +2.59%  "

$ws.Range("D32").Value = "'0.0610"
$ws.Range("E32").Value = "  -2.20%  "

$ws.Range("E33").Value = "  +2.33%  "

$ws.Range("D34").Value = "'0.0900"
$ws.Range("E34").Value = "  +2.03%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("E36").Value = "  -1.78%  "

$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("E38").Value = "  +5.09%  "

$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("D40").Value = "'3.16"
$ws.Range("E40").Value = "  +12.70%  "

$ws.Range("D41").Value = "'5.07"
$ws.Range("E41").Value = "  +27.56%  "

$ws.Range("E42").Value = "  -1.57%  "

$ws.Range("D43").Value = "'17.24"
$ws.Range("E43").Value = "  -5.50%  "

$ws.Range("E44").Value = "  -1.42%  "

$ws.Range("D45").Value = "'95.23"
$ws.Range("E45").Value = "  -0.86%  "

$ws.Range("D46").Value = "'2.42"
$ws.Range("E46").Value = "  +2.11%  "

$ws.Range("D47").Value = "1.274.92"
$ws.Range("E47").Value = "  -1.61%  "

$ws.Range("E48").Value = "  -2.29%  "

$ws.Range("D49").Value = "2.233.51"
$ws.Range("E49").Value = "  -0.44%  "

$ws.Range("E50").Value = "  -1.70%  "

$ws.Range("E51").Value = "  -18.09%  "
